# GANTT Chart update: convert the repeating "Week-1..Week-4" per-month header
# row into a single continuous week sequence (Week-1 .. Week-23), matching the
# commit "updated gant and added project milestones to proposal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New continuous week labels for the header row (row 2), columns C through Z
# (in order). Note column C and D both end up "Week-1", matching the source
# edit exactly.
$weekCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
$weekLabels = @("Week-1","Week-1","Week-2","Week-3","Week-4","Week-5","Week-6","Week-7","Week-8","Week-9","Week-10","Week-11","Week-12","Week-13","Week-14","Week-15","Week-16","Week-17","Week-18","Week-19","Week-20","Week-21","Week-22","Week-23")

for ($i = 0; $i -lt $weekCols.Length; $i++) {
    $ws.Range("$($weekCols[$i])2").Value = $weekLabels[$i]
}

# Column Z no longer ends a 4-week month block, so its right border changes
# from medium to thin (matching the other interior week-header cells), and
# the column widens slightly to fit the new "Week-23" text (stored OOXML
# width of 8.5 corresponds to a ColumnWidth of ~7.6667 in Excel's character
# units).
$ws.Range("Z2").Borders.Item(10).Weight = 2
$ws.Columns.Item(26).ColumnWidth = 7.6666667

# Update the selection to reflect where the author was working.
$ws.Range("D2:Z2").Select() | Out-Null
